$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: convert Age cell (B11) from text "21" to a real number 21
$ws.Range("B11").Value = 21

# New data rows to append (rows 12-17), Age stored as real numbers
$data = @(
    @("Solin Neat",     22, "Female", "Kompong Thom", "Class B 2025", "image\55b8f49f301c47a6b7ee3000cc00bcbe.png"),
    @("Bopha Khat",      20, "Female", "Kompot",        "Class B 2025", "image\fcfd397106ad4bbdbb133df0246a792d.png"),
    @("Bunyoung Hean",   23, "Male",   "Rotanakiri",    "Class A 2025", "image\778c5d5c26024057a335cd35381b9996.png"),
    @("Chhun Kimchhik",  20, "Female", "Bathombong",    "Class A 2025", "image\da72b4ddf65e47a8b18eaa5e866879c5.png"),
    @("Kanha Mao",       20, "Female", "Priveng ",      "Class B 2025", "image\ea1715a1a8354986b29cc68c6c0224d8.png"),
    @("Channak Kan",     20, "Female", "Bathombong ",   "Class B 2025", "image\ece17cec0596494f9747286afaab640d.png")
)

$r = 12
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

# Row 18: new last row, Age cell (B18) kept as text "21" (matches source quirk)
$ws.Cells.Item(18, 1).Value = "Hort Heng"
$ws.Cells.Item(18, 2).NumberFormat = "@"
$ws.Cells.Item(18, 2).Value = "21"
$ws.Cells.Item(18, 3).Value = "Male"
$ws.Cells.Item(18, 4).Value = "Priveng"
$ws.Cells.Item(18, 5).Value = "Class A 2025"
$ws.Cells.Item(18, 6).Value = "image\0909ecf105fe4905a3a7d5816b1bbc3b.png"
